# The commit removes the final slide of the deck (sldId 275 / rId20 /
# ppt/slides/slide19.xml, the "THANH VIEN NHOM 1 (TO 1)" member-list
# slide) from the slide list.
$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()
